$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Status column (E2:E11): "OPTIMAL" -> "TIME_LIMIT" (MP time limit change)
$ws.Range("E2").Value = "TIME_LIMIT"
$ws.Range("E3").Value = "TIME_LIMIT"
$ws.Range("E4").Value = "TIME_LIMIT"
$ws.Range("E5").Value = "TIME_LIMIT"
$ws.Range("E6").Value = "TIME_LIMIT"
$ws.Range("E7").Value = "TIME_LIMIT"
$ws.Range("E8").Value = "TIME_LIMIT"
$ws.Range("E9").Value = "TIME_LIMIT"
$ws.Range("E10").Value = "TIME_LIMIT"
$ws.Range("E11").Value = "TIME_LIMIT"

# Corrected fixed recourse data: objective (B), gap (C) and solve time (D)
$ws.Range("B2").Value = -876.2104421878598
$ws.Range("C2").Value = 17.5040401738034
$ws.Range("D2").Value = 6113.479721843

$ws.Range("B3").Value = -876.021618391906
$ws.Range("C3").Value = 16.387172082251393
$ws.Range("D3").Value = 6128.454741491

$ws.Range("B4").Value = -847.5899662464269
$ws.Range("C4").Value = 21.587613429401056
$ws.Range("D4").Value = 6093.922380239

$ws.Range("B5").Value = -866.6470038989896
$ws.Range("C5").Value = 18.7843993816022
$ws.Range("D5").Value = 6138.187756747

$ws.Range("B6").Value = -862.6312809711787
$ws.Range("C6").Value = 19.51898804737151
$ws.Range("D6").Value = 5891.220903613

$ws.Range("B7").Value = -869.3345895083158
$ws.Range("C7").Value = 18.64903338927495
$ws.Range("D7").Value = 5914.894766762

$ws.Range("B8").Value = -870.3357241901938
$ws.Range("C8").Value = 18.14478172461124
$ws.Range("D8").Value = 6099.909136346

$ws.Range("B9").Value = -851.210091750839
$ws.Range("C9").Value = 19.354555684943282
$ws.Range("D9").Value = 6108.129325466

$ws.Range("B10").Value = -853.1355957022353
$ws.Range("C10").Value = 20.570581618543063
$ws.Range("D10").Value = 6003.456472195

$ws.Range("B11").Value = -853.2554335268328
$ws.Range("C11").Value = 20.855492771536152
$ws.Range("D11").Value = 6077.215627943

$wb.Save()
